$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix spelling typos: paranthasis/paranthesis/parantheses -> parenthesis/parenthesis/parentheses
$ws.Range("C16").Value = "there must be a single space between the closing parenthesis and the opening brace of a multi-line IF statement found 0 spaces programming conventions"
$ws.Range("C25").Value = "closing parenthesis of a multi-line IF statement must be on a new line"
$ws.Range("C27").Value = "First condition of a multi-line IF statement must directly follow the opening parenthesis"
$ws.Range("C35").Value = "using statements but including a set of parentheses when they are not needed"
$ws.Range("C38").Value = "There must be a single space between the closing parenthesis and the opening brace of a multi-line IF statement found a new line"

# Add percentage column for the totals section
$ws.Range("C42").Formula = "=B42/`$B`$40"
$ws.Range("C43").Formula = "=B43/`$B`$40"
$ws.Range("C44").Formula = "=B44/`$B`$40"
$ws.Range("C45").Formula = "=B45/`$B`$40"
$ws.Range("C42:C45").Style = "Percent"
$ws.Range("C42:C45").NumberFormat = "0.00%"

$ws.Range("C47").Formula = "=SUM(C42:C45)"
$ws.Range("C47").NumberFormat = "0.00%"

# Update the view/selection
$ws.Application.ActiveWindow.ScrollRow = 38
$ws.Range("B47").Select()
